$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B59").Value = 286057
$ws.Range("E59").Value = 41735
$ws.Range("F59").Value = 30558
$ws.Range("G59").Value = 41906
$ws.Range("H59").Value = 32644
$ws.Range("I59").Value = 117559
$ws.Range("J59").Value = 59078
$ws.Range("P59").Value = 6351
$ws.Range("Q59").Value = 27930
$ws.Range("R59").Value = 226979
$ws.Range("U59").Value = 33128
$ws.Range("V59").Value = 26408
$ws.Range("W59").Value = 35018
$ws.Range("X59").Value = 26293
$ws.Range("Y59").Value = 89629
